$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O1").Value = "Multimedia Folder"
